$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 5).Value = 3
$ws.Cells.Item(2, 6).Value = 1
$ws.Cells.Item(2, 7).Value = 0.7304773333333333
$ws.Cells.Item(2, 8).Value = 2.191432
$ws.Cells.Item(2, 9).Value = 0.03163269997405359
$ws.Cells.Item(2, 10).Value = 0.03163269997405359
$ws.Cells.Item(2, 13).Value = 1.845768666666667
$ws.Cells.Item(2, 14).Value = 5.537306
$ws.Cells.Item(2, 15).Value = 0.01459089321241885
$ws.Cells.Item(2, 16).Value = 0.01459089321241885
$ws.Cells.Item(2, 17).Value = 1.348292173576889
$ws.Cells.Item(2, 18).Value = 12.134629562192
$ws.Cells.Item(2, 19).Value = 0.0004615493473419005
$ws.Cells.Item(2, 20).Value = 0.0004615493473419005
$ws.Cells.Item(3, 5).Value = 3
$ws.Cells.Item(3, 6).Value = 1
$ws.Cells.Item(3, 7).Value = 0.7304773333333333
$ws.Cells.Item(3, 8).Value = 2.191432
$ws.Cells.Item(3, 9).Value = 0.03163269997405359
$ws.Cells.Item(3, 10).Value = 0.03163269997405359
$ws.Cells.Item(3, 15).Value = 0.6557810310272387
$ws.Cells.Item(3, 16).Value = 0.6557810310272387
$ws.Cells.Item(3, 17).Value = 60.59837590762755
$ws.Cells.Item(3, 18).Value = 545.385383168648
$ws.Cells.Item(3, 19).Value = 0.02074412460316017
$ws.Cells.Item(3, 20).Value = 0.02074412460316017
$ws.Cells.Item(4, 5).Value = 3
$ws.Cells.Item(4, 6).Value = 1
$ws.Cells.Item(4, 7).Value = 0.7304773333333333
$ws.Cells.Item(4, 8).Value = 2.191432
$ws.Cells.Item(4, 9).Value = 0.03163269997405359
$ws.Cells.Item(4, 10).Value = 0.03163269997405359
$ws.Cells.Item(4, 13).Value = 41.69841866666667
$ws.Cells.Item(4, 14).Value = 125.095256
$ws.Cells.Item(4, 15).Value = 0.3296280757603424
$ws.Cells.Item(4, 16).Value = 0.3296280757603424
$ws.Cells.Item(4, 17).Value = 30.45974967184356
$ws.Cells.Item(4, 18).Value = 274.137747046592
$ws.Cells.Item(4, 19).Value = 0.01042702602355152
$ws.Cells.Item(4, 20).Value = 0.01042702602355152
$ws.Cells.Item(5, 9).Value = 0.4074771110502447
$ws.Cells.Item(5, 10).Value = 0.4074771110502448
$ws.Cells.Item(5, 13).Value = 1.845768666666667
$ws.Cells.Item(5, 14).Value = 5.537306
$ws.Cells.Item(5, 15).Value = 0.01459089321241885
$ws.Cells.Item(5, 16).Value = 0.01459089321241885
$ws.Cells.Item(5, 17).Value = 17.36804636314333
$ws.Cells.Item(5, 18).Value = 156.31241726829
$ws.Cells.Item(5, 19).Value = 0.005945455013839057
$ws.Cells.Item(5, 20).Value = 0.005945455013839058
$ws.Cells.Item(6, 9).Value = 0.4074771110502447
$ws.Cells.Item(6, 10).Value = 0.4074771110502448
$ws.Cells.Item(6, 15).Value = 0.6557810310272387
$ws.Cells.Item(6, 16).Value = 0.6557810310272387
$ws.Cells.Item(6, 19).Value = 0.2672157600045301
$ws.Cells.Item(6, 20).Value = 0.2672157600045302
$ws.Cells.Item(7, 9).Value = 0.4074771110502447
$ws.Cells.Item(7, 10).Value = 0.4074771110502448
$ws.Cells.Item(7, 13).Value = 41.69841866666667
$ws.Cells.Item(7, 14).Value = 125.095256
$ws.Cells.Item(7, 15).Value = 0.3296280757603424
$ws.Cells.Item(7, 16).Value = 0.3296280757603424
$ws.Cells.Item(7, 17).Value = 392.3677336988933
$ws.Cells.Item(7, 18).Value = 3531.30960329004
$ws.Cells.Item(7, 19).Value = 0.1343158960318755
$ws.Cells.Item(7, 20).Value = 0.1343158960318755
$ws.Cells.Item(8, 7).Value = 12.95234266666667
$ws.Cells.Item(8, 8).Value = 38.857028
$ws.Cells.Item(8, 9).Value = 0.5608901889757016
$ws.Cells.Item(8, 10).Value = 0.5608901889757018
$ws.Cells.Item(8, 13).Value = 1.845768666666667
$ws.Cells.Item(8, 14).Value = 5.537306
$ws.Cells.Item(8, 15).Value = 0.01459089321241885
$ws.Cells.Item(8, 16).Value = 0.01459089321241885
$ws.Cells.Item(8, 17).Value = 23.90702825406311
$ws.Cells.Item(8, 18).Value = 215.163254286568
$ws.Cells.Item(8, 19).Value = 0.008183888851237891
$ws.Cells.Item(8, 20).Value = 0.008183888851237893
$ws.Cells.Item(9, 7).Value = 12.95234266666667
$ws.Cells.Item(9, 8).Value = 38.857028
$ws.Cells.Item(9, 9).Value = 0.5608901889757016
$ws.Cells.Item(9, 10).Value = 0.5608901889757018
$ws.Cells.Item(9, 15).Value = 0.6557810310272387
$ws.Cells.Item(9, 16).Value = 0.6557810310272387
$ws.Cells.Item(9, 17).Value = 1074.490465320032
$ws.Cells.Item(9, 18).Value = 9670.414187880293
$ws.Cells.Item(9, 19).Value = 0.3678211464195484
$ws.Cells.Item(9, 20).Value = 0.3678211464195484
$ws.Cells.Item(10, 7).Value = 12.95234266666667
$ws.Cells.Item(10, 8).Value = 38.857028
$ws.Cells.Item(10, 9).Value = 0.5608901889757016
$ws.Cells.Item(10, 10).Value = 0.5608901889757018
$ws.Cells.Item(10, 13).Value = 41.69841866666667
$ws.Cells.Item(10, 14).Value = 125.095256
$ws.Cells.Item(10, 15).Value = 0.3296280757603424
$ws.Cells.Item(10, 16).Value = 0.3296280757603424
$ws.Cells.Item(10, 17).Value = 540.0922072287965
$ws.Cells.Item(10, 18).Value = 4860.829865059168
$ws.Cells.Item(10, 19).Value = 0.1848851537049153
$ws.Cells.Item(10, 20).Value = 0.1848851537049154
